# Applies the upstream "Add files via upload" update to the Championship
# 2023-2024 results sheet:
#   - Two pairs of match rows had their match-detail columns (F:V) swapped
#     back into the correct chronological order (the Indice/date columns
#     A:E were already correct and stay untouched).
#   - One new match row (Birmingham 0-0 QPR) is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match details (columns F:V) between row 18 and row 23 ---------
$row18 = $ws.Range("F18:V18").Value2
$row23 = $ws.Range("F23:V23").Value2
$ws.Range("F18:V18").Value2 = $row23
$ws.Range("F23:V23").Value2 = $row18

# --- Swap match details (columns F:V) between row 74 and row 77 ---------
$row74 = $ws.Range("F74:V74").Value2
$row77 = $ws.Range("F77:V77").Value2
$ws.Range("F74:V74").Value2 = $row77
$ws.Range("F77:V77").Value2 = $row74

# --- Append new row 86 (Birmingham 0-0 QPR) ------------------------------
# Copy the formatting of the previous last row (85) so the new row matches
# the existing styling (bold/bordered index column, date-formatted column E).
$ws.Range("A85:V85").Copy()
$ws.Range("A86:V86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A86").Value2 = 85
$ws.Range("B86").Value2 = "england"
$ws.Range("C86").Value2 = "championship"
$ws.Range("D86").Value2 = "2023-2024"
$ws.Range("E86").Value2 = 45191.875
$ws.Range("F86").Value2 = "Birmingham"
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = "QPR"
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 1.99
$ws.Range("K86").Value2 = "18/09/2023 19:42"
$ws.Range("L86").Value2 = 1.9
$ws.Range("M86").Value2 = "22/09/2023 20:55"
$ws.Range("N86").Value2 = 3.6
$ws.Range("O86").Value2 = "18/09/2023 19:42"
$ws.Range("P86").Value2 = 3.56
$ws.Range("Q86").Value2 = "22/09/2023 20:57"
$ws.Range("R86").Value2 = 3.92
$ws.Range("S86").Value2 = "18/09/2023 19:42"
$ws.Range("T86").Value2 = 4.48
$ws.Range("U86").Value2 = "22/09/2023 20:58"
$ws.Range("V86").Value2 = "https://www.betexplorer.com/football/england/championship/birmingham-qpr/6oZ4XXkU/"
